# Update the cryptos worksheet with latest scraped prices / volume percentages.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price/Volume columns to be treated as plain text while we write the
# new values, so strings like "324.02" or "1.038" are not silently reinterpreted
# as numbers (which would introduce floating point noise / strip the formatting).
$ws.Range("D2:E51").NumberFormat = "@"

# Row 2 - Bitcoin
$ws.Range("D2").Value = "27.689.33"
$ws.Range("E2").Value = "  +2.80%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "1.861.34"
$ws.Range("E3").Value = "  +2.74%  "

# Row 4 - TetherUSD
$ws.Range("D4").Value = "1.038"
$ws.Range("E4").Value = "  +2.92%  "

# Row 5 - BNB
$ws.Range("D5").Value = "324.02"
$ws.Range("E5").Value = "  +3.69%  "

# Row 6 - USDC
$ws.Range("D6").Value = "1.034"
$ws.Range("E6").Value = "  +2.80%  "

# Row 7 - XRP
$ws.Range("D7").Value = "0.4414"
$ws.Range("E7").Value = "  +2.93%  "

# Row 8 - Cardano
$ws.Range("D8").Value = "0.3806"
$ws.Range("E8").Value = "  +3.06%  "

# Row 9 - Dogecoin (only volume changes)
$ws.Range("E9").Value = "  +2.73%  "

# Row 10 - Polygon (only price changes)
$ws.Range("D10").Value = "0.8827"

# Row 11 - Solana
$ws.Range("D11").Value = "21.70"
$ws.Range("E11").Value = "  +2.45%  "

# Row 12 - WrappedEther
$ws.Range("D12").Value = "1.874.98"
$ws.Range("E12").Value = "  -7.63%  "

# Row 13 - Polkadot
$ws.Range("D13").Value = "5.552"
$ws.Range("E13").Value = "  +3.06%  "

# Row 14 - Chainlink
$ws.Range("D14").Value = "6.745"
$ws.Range("E14").Value = "  +1.75%  "

# Row 15 - TRON
$ws.Range("D15").Value = "0.07214"
$ws.Range("E15").Value = "  +4.50%  "

# Row 16 - Litecoin
$ws.Range("D16").Value = "83.66"
$ws.Range("E16").Value = "  +3.68%  "

# Row 17 - BinanceUSD
$ws.Range("D17").Value = "1.040"
$ws.Range("E17").Value = "  +3.50%  "

# Row 18 - ShibaInu
$ws.Range("D18").Value = "0.000009091"
$ws.Range("E18").Value = "  +1.96%  "

# Row 19 - Dai (only price changes)
$ws.Range("D19").Value = "1.034"

# Row 21 - WrappedBTC
$ws.Range("D21").Value = "27.725.29"
$ws.Range("E21").Value = "  +2.82%  "

# Row 22 - Uniswap
$ws.Range("D22").Value = "5.301"
$ws.Range("E22").Value = "  +2.01%  "

# Row 23 - Cosmos
$ws.Range("D23").Value = "11.42"
$ws.Range("E23").Value = "  +4.41%  "

# Row 24 - Monero
$ws.Range("D24").Value = "158.71"
$ws.Range("E24").Value = "  +3.23%  "

# Row 25 - Toncoin
$ws.Range("D25").Value = "1.933"
$ws.Range("E25").Value = "  +2.48%  "

# Row 26 - EthereumClassic
$ws.Range("D26").Value = "18.84"
$ws.Range("E26").Value = "  +2.82%  "

# Row 27 - LidoDAOToken
$ws.Range("D27").Value = "1.994"
$ws.Range("E27").Value = "  +5.14%  "

# Row 28 - InternetComputer(DFINITY)
$ws.Range("D28").Value = "5.313"
$ws.Range("E28").Value = "  +1.64%  "

# Row 29 - BitcoinCash
$ws.Range("D29").Value = "117.59"
$ws.Range("E29").Value = "  +2.15%  "

# Row 30 - Stellar
$ws.Range("D30").Value = "0.09084"
$ws.Range("E30").Value = "  +1.72%  "

# Row 31 - ARBITRUM
$ws.Range("D31").Value = "1.210"
$ws.Range("E31").Value = "  +4.58%  "

# Row 32 - ImmutableX
$ws.Range("D32").Value = "0.7657"
$ws.Range("E32").Value = "  +3.12%  "

# Row 33 - Filecoin
$ws.Range("D33").Value = "4.570"
$ws.Range("E33").Value = "  +3.39%  "

# Row 34 - HuobiToken (only volume changes)
$ws.Range("E34").Value = "  +3.40%  "

# Row 35 - Frax (only volume changes)
$ws.Range("E35").Value = "  +2.67%  "

# Row 36 - TrustWalletToken
$ws.Range("D36").Value = "1.160"
$ws.Range("E36").Value = "  +3.25%  "

# Row 37 - VeChain
$ws.Range("D37").Value = "0.01984"
$ws.Range("E37").Value = "  +3.38%  "

# Row 38 - Hedera (only volume changes)
$ws.Range("E38").Value = "  +2.59%  "

# Row 39 - was MXToken, now TheSandbox
$ws.Range("B39").Value = "TheSandbox"
$ws.Range("C39").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D39").Value = "0.5192"
$ws.Range("E39").Value = "  +2.11%  "

# Row 40 - was TheSandbox, now MXToken
$ws.Range("B40").Value = "MXToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D40").Value = "2.836"
$ws.Range("E40").Value = "  +3.73%  "

# Row 41 - Algorand
$ws.Range("D41").Value = "0.1689"
$ws.Range("E41").Value = "  +2.71%  "

# Row 42 - FraxShare
$ws.Range("D42").Value = "6.848"
$ws.Range("E42").Value = "  +6.59%  "

# Row 43 - Aptos
$ws.Range("D43").Value = "8.684"
$ws.Range("E43").Value = "  +5.31%  "

# Row 44 - Quant
$ws.Range("D44").Value = "109.51"
$ws.Range("E44").Value = "  +2.41%  "

# Row 45 - EnergySwap
$ws.Range("D45").Value = "10.62"
$ws.Range("E45").Value = "  +2.52%  "

# Row 46 - NEARProtocol (only volume changes)
$ws.Range("E46").Value = "  +4.45%  "

# Row 47 - Decentraland
$ws.Range("D47").Value = "0.4680"
$ws.Range("E47").Value = "  +2.81%  "

# Row 48 - Cronos (only volume changes)
$ws.Range("E48").Value = "  +2.21%  "

# Row 49 - RenderToken
$ws.Range("D49").Value = "1.862"
$ws.Range("E49").Value = "  +3.61%  "

# Row 50 - Elrond
$ws.Range("D50").Value = "39.68"
$ws.Range("E50").Value = "  +4.93%  "

# Row 51 - was ThetaToken, now Aave
$ws.Range("B51").Value = "Aave"
$ws.Range("C51").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D51").Value = "64.38"
$ws.Range("E51").Value = "  +1.33%  "

# Restore the default (unstyled) cell style now that the text values are set,
# matching the original workbook's formatting for these cells.
$ws.Range("D2:E51").Style = "Normal"
